$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Create the new "admonition" styles (Caution / Important / Note /
#    Tip / Warning) together with their linked character styles, mirroring
#    the relationships baked into the target styles.xml:
#      Caution    (paragraph) -> basedOn Abstract,  link CautionChar
#      CautionChar(character) -> basedOn AbstractChar, link Caution
#      Important  (paragraph) -> basedOn Caution,   link ImportantChar
#      ImportantChar(char)    -> basedOn CautionChar, link Important
#      Note       (paragraph) -> basedOn Caution,   link NoteChar
#      NoteChar(char)         -> basedOn CautionChar, link Note
#      Tip        (paragraph) -> basedOn Caution,   link TipChar
#      TipChar(char)          -> basedOn CautionChar, link Tip
#      Warning    (paragraph) -> basedOn Caution,   link WarningChar
#      WarningChar(char)      -> basedOn CautionChar, link Warning
# ---------------------------------------------------------------------

$caution = $d.Styles.Add("Caution", 1)
$caution.BaseStyle = $d.Styles.Item("Abstract")
$caution.QuickStyle = $true

$cautionChar = $d.Styles.Add("Caution Char", 2)
$cautionChar.BaseStyle = $d.Styles.Item("AbstractChar")
$caution.LinkStyle = $cautionChar
$cautionChar.LinkStyle = $caution

$important = $d.Styles.Add("Important", 1)
$important.BaseStyle = $caution
$important.QuickStyle = $true

$importantChar = $d.Styles.Add("Important Char", 2)
$importantChar.BaseStyle = $cautionChar
$important.LinkStyle = $importantChar
$importantChar.LinkStyle = $important

$note = $d.Styles.Add("Note", 1)
$note.BaseStyle = $caution
$note.QuickStyle = $true

$noteChar = $d.Styles.Add("Note Char", 2)
$noteChar.BaseStyle = $cautionChar
$note.LinkStyle = $noteChar
$noteChar.LinkStyle = $note

$tip = $d.Styles.Add("Tip", 1)
$tip.BaseStyle = $caution
$tip.QuickStyle = $true

$tipChar = $d.Styles.Add("Tip Char", 2)
$tipChar.BaseStyle = $cautionChar
$tip.LinkStyle = $tipChar
$tipChar.LinkStyle = $tip

$warning = $d.Styles.Add("Warning", 1)
$warning.BaseStyle = $caution
$warning.QuickStyle = $true

$warningChar = $d.Styles.Add("Warning Char", 2)
$warningChar.BaseStyle = $cautionChar
$warning.LinkStyle = $warningChar
$warningChar.LinkStyle = $warning

# ---------------------------------------------------------------------
# 2. Turn the trailing empty paragraph into "Admonitions: " and append the
#    five demo paragraphs (one per new style) right after it.
# ---------------------------------------------------------------------

$count = $d.Paragraphs.Count
$admon = $d.Paragraphs.Item($count)
$admon.Range.Text = "Admonitions: "

$pCaution = $d.Paragraphs.Add()
$pCaution.Range.Text = "Caution"
$pCaution.Range.Style = "Caution"

$pImportant = $d.Paragraphs.Add()
$pImportant.Range.Text = "Important"
$pImportant.Range.Style = "Important"

$pNote = $d.Paragraphs.Add()
$pNote.Range.Text = "Note"
$pNote.Range.Style = "Note"

$pTip = $d.Paragraphs.Add()
$pTip.Range.Text = "Tip"
$pTip.Range.Style = "Tip"

$pWarning = $d.Paragraphs.Add()
$pWarning.Range.Text = "Warning"
$pWarning.Range.Style = "Warning"
